# Port battle calculator correction.
#
# 1) "Deep water port": the SUM ranges in D3/E3 were missing the already
#    present row 35, extend them to match the sheet's real data range.
# 2) "Shallow water port": same fix for D3/E3 (missing row 21), plus the
#    ship list in column B (rows 4-21) was reshuffled/corrected and the BR
#    values (column C) for the affected ships were corrected to match.

$wb = $excel.ActiveWorkbook

# --- Deep water port: extend SUM ranges to include row 35 -----------------
$ws1 = $wb.Worksheets.Item("Deep water port")
$ws1.Range("D3").Formula = "=SUM(D4:D35)"
$ws1.Range("E3").Formula = "=SUM(E4:E35)"

# --- Shallow water port ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("Shallow water port")

# extend SUM ranges to include row 21
$ws2.Range("D3").Formula = "=SUM(D4:D21)"
$ws2.Range("E3").Formula = "=SUM(E4:E21)"

# corrected ship order / names for rows 4-21
$ws2.Range("B4").Value = "Hercules"
$ws2.Range("B5").Value = "Pandora"
$ws2.Range("B6").Value = "Mercury"
$ws2.Range("B7").Value = "Mortar Brig"
$ws2.Range("B8").Value = "NavyBrig"
$ws2.Range("B9").Value = "Niagara"
$ws2.Range("B10").Value = "Prince de Neufchatel"
$ws2.Range("B11").Value = "Rattlesnake"
$ws2.Range("B12").Value = "Rattlesnake Heavy"
$ws2.Range("B13").Value = "Snow"
$ws2.Range("B14").Value = "Brig"
$ws2.Range("B15").Value = "Pickle"
$ws2.Range("B16").Value = "Cutter"
$ws2.Range("B17").Value = "GunBoat"
$ws2.Range("B18").Value = "Lynx"
$ws2.Range("B19").Value = "Privateer"
$ws2.Range("B20").Value = "Yacht"
$ws2.Range("B21").Value = "Yacht Silver"

# corrected BR values for the ships whose rating changed
$ws2.Range("C4").Value = 100
$ws2.Range("C5").Value = 100
$ws2.Range("C13").Value = 80
$ws2.Range("C14").Value = 70
$ws2.Range("C15").Value = 55
$ws2.Range("C21").Value = 50
